$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4-22 down to 5-23.
# This also carries forward the date-column (D) cell style to the new row.
$ws.Range("A4:T4").EntireRow.Insert()

# Populate the newly-inserted row 4 with the new weekly data point.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44921
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103002
$ws.Range("J4").Value = "Ciruela"
$ws.Range("K4").Value = "Angeleno"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 450
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19111
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1062
$ws.Range("T4").Value = 18
